$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.000009858192104410727
$ws.Range("E2").Value = 0.000009858192104410727

$ws.Range("D3").Value = 0.745663216312449
$ws.Range("E3").Value = 0.745663216312449

$ws.Range("D4").Value = 0.002926450150719163
$ws.Range("E4").Value = 0.002926450150719163

$ws.Range("D5").Value = 0.00000000001117550847086308
$ws.Range("E5").Value = 0.00000000001117550847086308

$ws.Range("D6").Value = 0.4030865282386907
$ws.Range("E6").Value = 0.4030865282386907

$ws.Range("D7").Value = 0.8322219847239448
$ws.Range("E7").Value = 0.1677780152760552

$ws.Range("C8").Value = $False
$ws.Range("D8").Value = 0.000003249623138112965
$ws.Range("E8").Value = 0.9999967503768619

$ws.Range("C9").Value = $False
$ws.Range("D9").Value = 0.00000002234674711920502
$ws.Range("E9").Value = 0.9999999776532529

$ws.Range("C10").Value = $False
$ws.Range("D10").Value = 0.000002147287265778115
$ws.Range("E10").Value = 0.9999978527127342

$ws.Range("C11").Value = $False
$ws.Range("D11").Value = 0.2254630970790973
$ws.Range("E11").Value = 0.7745369029209027
$ws.Range("F11").Value = 4.6866135597229
$ws.Range("G11").Value = 0.5

$ws.Range("D12").Value = 0.0000000470266112855119
$ws.Range("E12").Value = 0.0000000470266112855119

$ws.Range("D13").Value = 0.9890497261575523
$ws.Range("E13").Value = 0.9890497261575523

$ws.Range("D14").Value = 0.003167714919631714
$ws.Range("E14").Value = 0.003167714919631714

$ws.Range("D15").Value = 0.00000000000000009661978871129893
$ws.Range("E15").Value = 0.00000000000000009661978871129893

$ws.Range("D16").Value = 0.1058022678936759
$ws.Range("E16").Value = 0.1058022678936759

$ws.Range("D17").Value = 0.8871086094102821
$ws.Range("E17").Value = 0.1128913905897179

$ws.Range("C18").Value = $False
$ws.Range("D18").Value = 0.0000000003524310740347391
$ws.Range("E18").Value = 0.9999999996475689

$ws.Range("C19").Value = $False
$ws.Range("D19").Value = 0.00000000002114786682176066
$ws.Range("E19").Value = 0.9999999999788521

$ws.Range("C20").Value = $False
$ws.Range("D20").Value = 0.000000001353072377516494
$ws.Range("E20").Value = 0.9999999986469276

$ws.Range("C21").Value = $False
$ws.Range("D21").Value = 0.04889888223341513
$ws.Range("E21").Value = 0.9511011177665849
$ws.Range("F21").Value = 7.453371524810791
$ws.Range("G21").Value = 0.5
